$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Date update
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 becomes Jurisdiction / United States of America
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the old duplicate "Contact" row (row 11); rows below shift up
$ws.Range("A11").EntireRow.Delete()

# Elements sheet: row 2 (top-level Extension) Short/Definition now use the
# real extension title/description instead of generic placeholders
$ws2.Range("K2").Value = "Coverage Days"
$ws2.Range("L2").Value = "Number of covered days of eligibility"
